# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold prices formatted as plain text (thousand-dot separators,
# leading/trailing zeros that must be preserved exactly). Writing a leading apostrophe
# forces Excel to store the value as text instead of coercing it to a Number, and
# resetting the style back to 'Normal' afterwards avoids leaving a stray NumberFormat
# style on the cell (keeps the cell style identical to the original file).
function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "29.750.01"
$ws.Range("E2").Value = "  +0.78%  "
Set-TextCell "D3" "1.851.86"
$ws.Range("E3").Value = "  +0.57%  "
Set-TextCell "D4" "0.9998"
$ws.Range("E4").Value = "  +0.13%  "
Set-TextCell "D5" "243.90"
$ws.Range("E5").Value = "  -0.34%  "
Set-TextCell "D6" "0.6549"
$ws.Range("E6").Value = "  +4.05%  "
Set-TextCell "D7" "1.001"
$ws.Range("E7").Value = "  +0.10%  "
Set-TextCell "D8" "48.15"
$ws.Range("E8").Value = "  +3.96%  "
Set-TextCell "D9" "0.07500"
$ws.Range("E9").Value = "  +0.69%  "
Set-TextCell "D10" "0.2979"
$ws.Range("E10").Value = "  +0.30%  "
Set-TextCell "D11" "24.55"
$ws.Range("E11").Value = "  +3.55%  "
Set-TextCell "D12" "0.07637"
$ws.Range("E12").Value = "  -0.45%  "
Set-TextCell "D13" "1.852.70"
$ws.Range("E13").Value = "  +0.74%  "
Set-TextCell "D14" "5.054"
$ws.Range("E14").Value = "  +0.62%  "
Set-TextCell "D15" "0.6864"
$ws.Range("E15").Value = "  +0.96%  "
Set-TextCell "D16" "83.54"
$ws.Range("E16").Value = "  -0.67%  "
Set-TextCell "D17" "0.000009680"
$ws.Range("E17").Value = "  +3.53%  "
Set-TextCell "D18" "6.126"
$ws.Range("E18").Value = "  +2.47%  "
Set-TextCell "D19" "29.769.04"
$ws.Range("E19").Value = "  +0.94%  "
Set-TextCell "D20" "2.103.74"
$ws.Range("E20").Value = "  +1.16%  "
Set-TextCell "D21" "237.62"
$ws.Range("E21").Value = "  -0.06%  "
Set-TextCell "D22" "12.63"
$ws.Range("E22").Value = "  +0.32%  "
Set-TextCell "D23" "1.0000"
$ws.Range("E23").Value = "  +0.05%  "
Set-TextCell "D24" "7.707"
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("E25").Value = "  +0.00%  "
Set-TextCell "D26" "158.23"
$ws.Range("E26").Value = "  -0.55%  "
Set-TextCell "D27" "0.1428"
$ws.Range("E27").Value = "  +0.88%  "
Set-TextCell "D28" "8.539"
$ws.Range("E28").Value = "  +0.11%  "
Set-TextCell "D29" "17.84"
$ws.Range("E29").Value = "  +0.16%  "
Set-TextCell "D30" "0.06081"
$ws.Range("E30").Value = "  +0.35%  "
Set-TextCell "D31" "1.492"
$ws.Range("E31").Value = "  -0.42%  "
Set-TextCell "D32" "1.277"
$ws.Range("E32").Value = "  +2.57%  "
Set-TextCell "D33" "4.137"
$ws.Range("E33").Value = "  +0.47%  "
Set-TextCell "D34" "4.075"
$ws.Range("E34").Value = "  -1.09%  "
Set-TextCell "D35" "1.875"
$ws.Range("E35").Value = "  +0.14%  "
Set-TextCell "D36" "1.183"
$ws.Range("E36").Value = "  +3.26%  "
Set-TextCell "D37" "0.7267"
$ws.Range("E37").Value = "  -0.26%  "
Set-TextCell "D38" "2.604"
$ws.Range("E38").Value = "  -0.41%  "
Set-TextCell "D39" "2.803"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("E40").Value = "  +1.41%  "
Set-TextCell "D41" "1.199.55"
Set-TextCell "D42" "6.270"
$ws.Range("E42").Value = "  -0.34%  "
Set-TextCell "D43" "0.9103"
$ws.Range("E43").Value = "  -0.54%  "
Set-TextCell "D44" "1.000"
$ws.Range("E44").Value = "  -0.05%  "
Set-TextCell "D45" "2.016.84"
$ws.Range("E45").Value = "  +1.16%  "
Set-TextCell "D46" "101.14"
$ws.Range("E46").Value = "  -0.78%  "
Set-TextCell "D47" "66.73"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D48" "0.00000000123"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D49" "7.309"
$ws.Range("E49").Value = "  +9.45%  "
Set-TextCell "D50" "0.4059"
$ws.Range("E50").Value = "  -0.33%  "
Set-TextCell "D51" "9.157"
$ws.Range("E51").Value = "  -1.42%  "
